$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("D2").Value = 0.9999999934554245
$ws.Range("E2").Value = 0.9999999934554245

# Row 3
$ws.Range("D3").Value = 0.9999840597829652
$ws.Range("E3").Value = 0.9999840597829652

# Row 4
$ws.Range("D4").Value = 0.22021000001724
$ws.Range("E4").Value = 0.22021000001724

# Row 6
$ws.Range("D6").Value = 0.9995209817050612
$ws.Range("E6").Value = 0.9995209817050612

# Row 7
$ws.Range("D7").Value = 0.9999999999992284
$ws.Range("E7").Value = 0.0000000000007716050021144838

# Row 8
$ws.Range("D8").Value = 0.9920148768728613
$ws.Range("E8").Value = 0.007985123127138727

# Row 9
$ws.Range("C9").Value = $false
$ws.Range("D9").Value = 0.00985588110446288
$ws.Range("E9").Value = 0.9901441188955371

# Row 10
$ws.Range("D10").Value = 0.9779544071085774
$ws.Range("E10").Value = 0.02204559289142261

# Row 11
$ws.Range("C11").Value = $false
$ws.Range("D11").Value = 0.0000000005130014740879528
$ws.Range("E11").Value = 0.9999999994869986
$ws.Range("F11").Value = 9.929945945739746
$ws.Range("G11").Value = 0.4

# Row 13
$ws.Range("D13").Value = 0.9989056102864499
$ws.Range("E13").Value = 0.9989056102864499

# Row 14
$ws.Range("D14").Value = 0.3445887013055929
$ws.Range("E14").Value = 0.3445887013055929

# Row 16
$ws.Range("D16").Value = 0.9999819771802413
$ws.Range("E16").Value = 0.9999819771802413

# Row 17
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0

# Row 18
$ws.Range("D18").Value = 0.9996666675087397
$ws.Range("E18").Value = 0.0003333324912603297

# Row 19
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = 0.0004391268535465651
$ws.Range("E19").Value = 0.9995608731464535

# Row 20
$ws.Range("D20").Value = 0.9949017237359963
$ws.Range("E20").Value = 0.00509827626400372

# Row 21
$ws.Range("C21").Value = $false
$ws.Range("D21").Value = 0.0000000009116627185631092
$ws.Range("E21").Value = 0.9999999990883373
$ws.Range("F21").Value = 15.96742534637451
$ws.Range("G21").Value = 0.4

$wb.Save()
